$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the header text in C1 from "biosampleNumber" to "bioSampleNumber"
$ws.Range("C1").Value = "bioSampleNumber"

# Restore the active selection to C1 (as reflected in the saved file)
$ws.Range("C1").Select()
